$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: copy the formatting (styles) of row 2 (A2:D2 -> A5:D5) so the date/
# duration cells keep the same number-format styles as the rest of the table,
# then overwrite the values with the new entry's data.
$ws.Range("A2:D2").Copy($ws.Range("A5:D5"))
$ws.Range("A5").Value = "Se crearon el resto de pestañas de la web administradora."
$ws.Range("B5").Value = 42870.791666666664
$ws.Range("C5").Value = 42870.958333333336
$ws.Range("D5").Value = 0.16666666666666666

# Row 6: same idea, copying the formatting of row 3 (which already has the
# taller row height used for wrapped, two-line entries).
$ws.Range("A3:D3").Copy($ws.Range("A6:D6"))
$ws.Range("A6").Value = "Se añadio la pestaña de empleados a la web administradora y se creó el formulario de login."
$ws.Range("B6").Value = 42871.833333333336
$ws.Range("C6").Value = 42871.958333333336
$ws.Range("D6").Value = 0.125
$ws.Rows("6").RowHeight = 30

# Update the active selection to D7, matching the saved sheet view.
$ws.Range("D7").Select()
